$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing formatting / precision).
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D19", "D22", "D23", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the refreshed crypto price feed.
$ws.Range('D2').Value = '26.526.00'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.839.25'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '258.53'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '0.5219'
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').Value = '0.3166'
$ws.Range('E8').Value = '  -3.14%  '
$ws.Range('D9').Value = '0.06776'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').Value = '18.67'
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').Value = '0.7777'
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('D12').Value = '0.07764'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').Value = '1.827.94'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').Value = '87.72'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').Value = '5.003'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').Value = '13.83'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '0.000007914'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('D20').Value = '26.552.63'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').Value = '2.074.35'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').Value = '4.596'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').Value = '5.959'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '9.310'
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '1.672'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('D28').Value = '16.88'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').Value = '111.70'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').Value = '4.161'
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('D31').Value = '0.08720'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '4.060'
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('D33').Value = '0.04874'
$ws.Range('E33').Value = '  +1.19%  '
$ws.Range('D34').Value = '1.132'
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').Value = '0.7197'
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('D36').Value = '2.860'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('D37').Value = '3.091'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '2.224'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').Value = '0.01728'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('D40').Value = '0.4812'
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('D41').Value = '0.8956'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').Value = '110.25'
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('D43').Value = '5.921'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '7.620'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').Value = '0.4158'
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').Value = '8.935'
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05825'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.1229'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('D50').Value = '34.78'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('D51').Value = '0.8900'
$ws.Range('E51').Value = '  +0.14%  '

Write-Output "done"
